# Auto-generated: refresh cached market-data values (currentAveragePrice* / Leve price & profit columns)
# per scheduled-runner update. Each row below corresponds to one <c> value changed in the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3216.2307
$ws.Range("J17").Value = 2770.6667
$ws.Range("L17").Value = 8312.000100000001
$ws.Range("N17").Value = -8648.000100000001
$ws.Range("H19").Value = 1766.3125
$ws.Range("J19").Value = 2263.3635
$ws.Range("L19").Value = 2263.3635
$ws.Range("N19").Value = -2613.3635
$ws.Range("H74").Value = 4416.5
$ws.Range("I74").Value = 4299.8
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 4299.8
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -3363.8
$ws.Range("N74").Value = -6872
$ws.Range("H77").Value = 4416.5
$ws.Range("I77").Value = 4299.8
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 21499
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -16819
$ws.Range("N77").Value = -34360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7014.3193
$ws.Range("I32").Value = 5002.2905
$ws.Range("K32").Value = 5002.2905
$ws.Range("M32").Value = -4715.2905
$ws.Range("H45").Value = 1156.7407
$ws.Range("I45").Value = 915.8461
$ws.Range("J45").Value = 1380.4286
$ws.Range("K45").Value = 915.8461
$ws.Range("L45").Value = 1380.4286
$ws.Range("M45").Value = -538.8461
$ws.Range("N45").Value = -2134.4286
$ws.Range("H61").Value = 7660.5
$ws.Range("I61").Value = 5999.85
$ws.Range("K61").Value = 5999.85
$ws.Range("M61").Value = -5787.85
$ws.Range("H122").Value = 1355.6428
$ws.Range("J122").Value = 2750
$ws.Range("L122").Value = 8250
$ws.Range("N122").Value = -13150
$ws.Range("H136").Value = 7660.5
$ws.Range("I136").Value = 5999.85
$ws.Range("K136").Value = 17999.55
$ws.Range("M136").Value = -15449.55

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 26981.23
$ws.Range("I82").Value = 13814
$ws.Range("J82").Value = 32833.332
$ws.Range("K82").Value = 13814
$ws.Range("L82").Value = 32833.332
$ws.Range("M82").Value = -13431
$ws.Range("N82").Value = -33599.332
$ws.Range("H85").Value = 26981.23
$ws.Range("I85").Value = 13814
$ws.Range("J85").Value = 32833.332
$ws.Range("K85").Value = 13814
$ws.Range("L85").Value = 32833.332
$ws.Range("M85").Value = -12488
$ws.Range("N85").Value = -35485.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1742.3125
$ws.Range("I31").Value = 1126.8
$ws.Range("J31").Value = 2768.1667
$ws.Range("K31").Value = 1126.8
$ws.Range("L31").Value = 2768.1667
$ws.Range("M31").Value = -831.8
$ws.Range("N31").Value = -3358.1667
$ws.Range("H34").Value = 1742.3125
$ws.Range("I34").Value = 1126.8
$ws.Range("J34").Value = 2768.1667
$ws.Range("K34").Value = 1126.8
$ws.Range("L34").Value = 2768.1667
$ws.Range("M34").Value = -924.8
$ws.Range("N34").Value = -3172.1667
$ws.Range("H86").Value = 2869.4285
$ws.Range("I86").Value = 2571.75
$ws.Range("K86").Value = 2571.75
$ws.Range("M86").Value = -1448.75
$ws.Range("H89").Value = 2869.4285
$ws.Range("I89").Value = 2571.75
$ws.Range("K89").Value = 12858.75
$ws.Range("M89").Value = -7242.75
$ws.Range("H99").Value = 2872.5
$ws.Range("I99").Value = 2663.3333
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 2663.3333
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = -1165.3333
$ws.Range("N99").Value = -6496
$ws.Range("H107").Value = 762.13336
$ws.Range("I107").Value = 768.2308
$ws.Range("K107").Value = 768.2308
$ws.Range("M107").Value = 1151.7692
$ws.Range("H122").Value = 3406.6
$ws.Range("I122").Value = 2185.25
$ws.Range("J122").Value = 4220.8335
$ws.Range("K122").Value = 6555.75
$ws.Range("L122").Value = 12662.5005
$ws.Range("M122").Value = -4105.75
$ws.Range("N122").Value = -17562.5005
$ws.Range("H126").Value = 2872.5
$ws.Range("I126").Value = 2663.3333
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 7989.999899999999
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -5519.999899999999
$ws.Range("N126").Value = -15440
$ws.Range("H134").Value = 2260.9756
$ws.Range("I134").Value = 2113.6758
$ws.Range("K134").Value = 6341.0274
$ws.Range("M134").Value = -3806.0274

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2245.389
$ws.Range("I140").Value = 1401.4
$ws.Range("K140").Value = 4204.200000000001
$ws.Range("M140").Value = 975.7999999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2080.8572
$ws.Range("I102").Value = 2135.1765
$ws.Range("K102").Value = 2135.1765
$ws.Range("M102").Value = -513.1765

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1251.5
$ws.Range("I22").Value = 683.6667
$ws.Range("J22").Value = 1716.091
$ws.Range("K22").Value = 683.6667
$ws.Range("L22").Value = 1716.091
$ws.Range("M22").Value = -388.6667
$ws.Range("N22").Value = -2306.091
$ws.Range("H27").Value = 1251.5
$ws.Range("I27").Value = 683.6667
$ws.Range("J27").Value = 1716.091
$ws.Range("K27").Value = 683.6667
$ws.Range("L27").Value = 1716.091
$ws.Range("M27").Value = -576.6667
$ws.Range("N27").Value = -1930.091
$ws.Range("H40").Value = 18609.666
$ws.Range("I40").Value = 25300.8
$ws.Range("K40").Value = 25300.8
$ws.Range("M40").Value = -25164.8
$ws.Range("H61").Value = 3906.1904
$ws.Range("I61").Value = 3542.7693
$ws.Range("K61").Value = 3542.7693
$ws.Range("M61").Value = -3340.7693
$ws.Range("H68").Value = 2041.75
$ws.Range("I68").Value = 1513
$ws.Range("J68").Value = 2782
$ws.Range("K68").Value = 1513
$ws.Range("L68").Value = 2782
$ws.Range("M68").Value = -764
$ws.Range("N68").Value = -4280
$ws.Range("H71").Value = 2041.75
$ws.Range("I71").Value = 1513
$ws.Range("J71").Value = 2782
$ws.Range("K71").Value = 7565
$ws.Range("L71").Value = 13910
$ws.Range("M71").Value = -3821
$ws.Range("N71").Value = -21398
$ws.Range("H82").Value = 2873.111
$ws.Range("I82").Value = 2000
$ws.Range("K82").Value = 2000
$ws.Range("M82").Value = -1639
$ws.Range("H85").Value = 2873.111
$ws.Range("I85").Value = 2000
$ws.Range("K85").Value = 2000
$ws.Range("M85").Value = -752
$ws.Range("H113").Value = 3906.1904
$ws.Range("I113").Value = 3542.7693
$ws.Range("K113").Value = 3542.7693
$ws.Range("M113").Value = -1372.7693
$ws.Range("H122").Value = 2340.4
$ws.Range("I122").Value = 2313
$ws.Range("J122").Value = 2450
$ws.Range("K122").Value = 6939
$ws.Range("L122").Value = 7350
$ws.Range("M122").Value = -4489
$ws.Range("N122").Value = -12250
$ws.Range("H136").Value = 1986.5938
$ws.Range("I136").Value = 1187.8636
$ws.Range("J136").Value = 3743.8
$ws.Range("K136").Value = 3563.5908
$ws.Range("L136").Value = 11231.4
$ws.Range("M136").Value = -1013.5908
$ws.Range("N136").Value = -16331.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 662.125
$ws.Range("I107").Value = 554.6
$ws.Range("K107").Value = 1663.8
$ws.Range("M107").Value = 256.1999999999998
$ws.Range("H122").Value = 38611.43
$ws.Range("I122").Value = 47212.06
$ws.Range("J122").Value = 2058.75
$ws.Range("K122").Value = 141636.18
$ws.Range("L122").Value = 6176.25
$ws.Range("M122").Value = -139186.18
$ws.Range("N122").Value = -11076.25
